# Updates the D (Price) and E (Volume 1h) text columns for the cryptos sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.403.67"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.773.77"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "354.17"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.81"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.550"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.13%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.586"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.62"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.95%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.01"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.28%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.19%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.211.03"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.780.56"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.922"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.410.56"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.61"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.59%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.15"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0962"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "265.60"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.12%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.99"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.161"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +12.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.21"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.69"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +8.45%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.14"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +7.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "51.78"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.68%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.51"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0824"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.31%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.22%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.42%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.53"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.73"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.96"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.096.50"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.24"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.901"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.79%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +8.52%  "
